# Auto-generated Excel COM-interop script
# Updates market-price-derived columns (H-N) on several sheets
# per the scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4252.6313
$ws.Range("I40").Value = 3300
$ws.Range("J40").Value = 4364.706
$ws.Range("K40").Value = 3300
$ws.Range("L40").Value = 4364.706
$ws.Range("M40").Value = -3125
$ws.Range("N40").Value = -4714.706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 696
$ws.Range("I4").Value = 235.6
$ws.Range("K4").Value = 235.6
$ws.Range("M4").Value = -119.6

$ws.Range("H61").Value = 21788930
$ws.Range("I61").Value = 38465464
$ws.Range("K61").Value = 38465464
$ws.Range("M61").Value = -38465252

$ws.Range("H82").Value = 51944.5
$ws.Range("J82").Value = 51944.5
$ws.Range("L82").Value = 51944.5
$ws.Range("N82").Value = -52666.5

$ws.Range("H85").Value = 51944.5
$ws.Range("J85").Value = 51944.5
$ws.Range("L85").Value = 51944.5
$ws.Range("N85").Value = -54440.5

$ws.Range("H136").Value = 21788930
$ws.Range("I136").Value = 38465464
$ws.Range("K136").Value = 115396392
$ws.Range("M136").Value = -115393842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2605.5454
$ws.Range("I99").Value = 2095.6667
$ws.Range("J99").Value = 4900
$ws.Range("K99").Value = 2095.6667
$ws.Range("L99").Value = 4900
$ws.Range("M99").Value = -597.6667000000002
$ws.Range("N99").Value = -7896

$ws.Range("H134").Value = 202002.4
$ws.Range("I134").Value = 2499.5
$ws.Range("K134").Value = 7498.5
$ws.Range("M134").Value = -4963.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 64940
$ws.Range("J87").Value = 64880
$ws.Range("L87").Value = 64880
$ws.Range("N87").Value = -67252

$ws.Range("H88").Value = 30124.75
$ws.Range("J88").Value = 30124.75
$ws.Range("L88").Value = 30124.75
$ws.Range("N88").Value = -30936.75

$ws.Range("H90").Value = 64940
$ws.Range("J90").Value = 64880
$ws.Range("L90").Value = 194640
$ws.Range("N90").Value = -206496

$ws.Range("H91").Value = 30124.75
$ws.Range("J91").Value = 30124.75
$ws.Range("L91").Value = 30124.75
$ws.Range("N91").Value = -32932.75

$ws.Range("H108").Value = 76718.39999999999
$ws.Range("J108").Value = 76718.39999999999
$ws.Range("L108").Value = 76718.39999999999
$ws.Range("N108").Value = -84398.39999999999

$ws.Range("H111").Value = 55000
$ws.Range("J111").Value = 55000
$ws.Range("L111").Value = 55000
$ws.Range("N111").Value = -63180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 4595.2
$ws.Range("I44").Value = 4595.2
$ws.Range("K44").Value = 13785.6
$ws.Range("M44").Value = -13387.6

$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -8189
$ws.Range("N69").Value = -16622

$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -22944
$ws.Range("N72").Value = -53112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 270000
$ws.Range("J39").Value = 40000
$ws.Range("L39").Value = 40000
$ws.Range("N39").Value = -41064

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null

$ws.Range("H132").Value = 66669588
$ws.Range("I132").Value = 71431490
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 214294470
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -214291940
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2577.125
$ws.Range("I16").Value = 2603.1667
$ws.Range("J16").Value = 2499
$ws.Range("K16").Value = 2603.1667
$ws.Range("L16").Value = 2499
$ws.Range("M16").Value = -2433.1667
$ws.Range("N16").Value = -2839

$ws.Range("H22").Value = 2542.7144
$ws.Range("I22").Value = 3240
$ws.Range("J22").Value = 799.5
$ws.Range("K22").Value = 3240
$ws.Range("L22").Value = 799.5
$ws.Range("M22").Value = -2945
$ws.Range("N22").Value = -1389.5

$ws.Range("H27").Value = 2542.7144
$ws.Range("I27").Value = 3240
$ws.Range("J27").Value = 799.5
$ws.Range("K27").Value = 3240
$ws.Range("L27").Value = 799.5
$ws.Range("M27").Value = -3133
$ws.Range("N27").Value = -1013.5

$ws.Range("H122").Value = 6493.375
$ws.Range("I122").Value = 5639.4
$ws.Range("J122").Value = 7916.6665
$ws.Range("K122").Value = 16918.2
$ws.Range("L122").Value = 23749.9995
$ws.Range("M122").Value = -14468.2
$ws.Range("N122").Value = -28649.9995

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws.Range("H132").Value = 74136.17999999999
$ws.Range("I132").Value = 47809.453
$ws.Range("K132").Value = 143428.359
$ws.Range("M132").Value = -140898.359

$ws.Range("H136").Value = 85136.234
$ws.Range("I136").Value = 50263.715
$ws.Range("K136").Value = 150791.145
$ws.Range("M136").Value = -148241.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null

$ws.Range("H64").Value = 64997
$ws.Range("J64").Value = 64997
$ws.Range("L64").Value = 64997
$ws.Range("N64").Value = -65493

$ws.Range("H67").Value = 64997
$ws.Range("J67").Value = 64997
$ws.Range("L67").Value = 64997
$ws.Range("N67").Value = -66713

$ws.Range("H93").Value = 77857
$ws.Range("J93").Value = 77857
$ws.Range("L93").Value = 77857
$ws.Range("N93").Value = -82849

$ws.Range("H112").Value = 99000
$ws.Range("J112").Value = 99000
$ws.Range("L112").Value = 99000
$ws.Range("N112").Value = -101954

$ws.Range("H136").Value = 52002
$ws.Range("I136").Value = 52002
$ws.Range("K136").Value = 156006
$ws.Range("M136").Value = -153456
